$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in L1 from "Estatura" to "Altura_cm"
$ws.Range("L1").Value = "Altura_cm"

# Fill in the new height data for rows 3-7
$ws.Range("L3").Value = 170
$ws.Range("L4").Value = 176
$ws.Range("L5").Value = 160
$ws.Range("L6").Value = 160
$ws.Range("L7").Value = 165

# Autofit the new column so its width matches data (bestFit)
$ws.Range("L1:L15").EntireColumn.AutoFit()

# Update selection to L8
$ws.Range("L8").Select()
